# Update the "想去人数" (F column) figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 2203
    5  = 13150
    7  = 116
    8  = 516
    11 = 987
    13 = 14372
    15 = 172
    21 = 36
    25 = 5428
    26 = 938
    27 = 22
    29 = 21
    30 = 53
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
